# CI: Auto Update Data (#101)
# Update MAA operator data cells with refreshed statistics and refresh timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "maa://24702 (94.18), maa://25390 (95.83), maa://36681 (86.3)"
$ws.Range("T3").Value = "maa://24617 (88.46), **maa://20790 (43.94), ***maa://37170 (20.0)"
$ws.Range("AF4").Value = "*maa://30062 (60.0), ***maa://26209 (13.04), *maa://39394 (66.67)"
$ws.Range("A8").Value = "更新日期：2024.11.20 13:18:41"
$ws.Range("X8").Value = "maa://21411 (95.92)"
$ws.Range("D13").Value = "maa://24999 (91.61), maa://36673 (92.42), maa://25001 (85.51)"
$ws.Range("P13").Value = "maa://22676 (91.59), *maa://22583 (75.41), *maa://22500 (56.82)"
$ws.Range("AF13").Value = "**maa://22737 (30.37), maa://39883 (91.67), *maa://39885 (56.0)"
$ws.Range("L14").Value = "maa://26245 (96.27), maa://21288 (96.21), maa://36682 (97.3), maa://39841 (94.03)"
$ws.Range("D15").Value = "*maa://22743 (77.13), maa://22734 (83.76), *maa://30808 (63.93), ***maa://36048 (25.0)"
$ws.Range("AF15").Value = "maa://21364 (80.54), *maa://22766 (70.37), *maa://36666 (78.21)"
$ws.Range("D16").Value = "maa://21441 (96.26), maa://36679 (92.68), maa://37650 (96.77)"
$ws.Range("T16").Value = "maa://22729 (95.3), *maa://28648 (67.24), maa://36674 (84.21)"
$ws.Range("H17").Value = "maa://22430 (88.46), maa://39599 (84.38)"
$ws.Range("D18").Value = "maa://24570 (96.95)"
$ws.Range("AB19").Value = "*maa://30709 (62.56), *maa://36668 (53.42)"
$ws.Range("L23").Value = "maa://39756 (93.02), maa://39875 (93.22)"
$ws.Range("X28").Value = "maa://39929 (89.2), ***maa://39723 (14.29), maa://41749 (85.71)"
$ws.Range("L29").Value = "maa://28432 (93.38), *maa://28440 (73.81), maa://31400 (100.0), *maa://28650 (71.43)"
$ws.Range("AB30").Value = "maa://42979 (97.06)"
$ws.Range("H32").Value = "maa://21895 (97.09), maa://36667 (98.28), **maa://20793 (38.78), maa://22760 (100.0)"
$ws.Range("T32").Value = "maa://41108 (87.5), maa://42859 (93.62), maa://41238 (94.92)"
$ws.Range("AF38").Value = "maa://36697 (85.62)"
$ws.Range("P39").Value = "maa://24709 (91.38)"
$ws.Range("H46").Value = "maa://35931 (92.42)"
$ws.Range("H53").Value = "maa://32534 (93.31), **maa://32434 (34.78)"
$ws.Range("H55").Value = "maa://32532 (92.31)"
